# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value mapping (column F)
$updates = @{
    2  = 71
    3  = 1059
    5  = 3053
    7  = 2182
    8  = 187
    9  = 108
    10 = 1029
    12 = 37
    13 = 252
    14 = 93
    15 = 96
    16 = 43
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
